# Apply edit: "Add classHoSoNV, classChucVu.Test Mapping to SQL"
#
# 1) In the "Chức vụ" class block, the "(khóa chính)" annotation moves
#    from the -MaChucVu line to the -TenChucVu line:
#       -MaChucVu (khóa chính)  ->  -MaChucVu
#       -TenChucVu               ->  -TenChucVu (khóa chính)
#
# 2) In the "Bảng Lương" class block, the "(khóa ngoại)" annotation is
#    removed from the -HeSoLuong line (leaving a trailing space):
#       -HeSoLuong (khóa ngoại)  ->  -HeSoLuong

$d = $word.ActiveDocument

function Get-ParaPlainText($p) {
    $t = $p.Range.Text
    # Drop the trailing paragraph-mark character (CR, code 13) so we can
    # compare the paragraph's visible text exactly.
    if ($t.Length -gt 0 -and [int][char]$t[$t.Length - 1] -eq 13) {
        $t = $t.Substring(0, $t.Length - 1)
    }
    return $t
}

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = Get-ParaPlainText $p

    if ($txt -eq "-MaChucVu (khóa chính)") {
        $p.Range.Text = "-MaChucVu "
    }
    elseif ($txt -eq "-TenChucVu") {
        $p.Range.Text = "-TenChucVu (khóa chính)"
    }
    elseif ($txt -eq "-HeSoLuong (khóa ngoại)") {
        $p.Range.Text = "-HeSoLuong "
    }
}
